$wb = $excel.ActiveWorkbook

# --- Sheet2: update release dates (column E) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("E3").Value = "'12/3/2022"
$ws2.Range("E3").NumberFormat = "mm-dd-yy"

$ws2.Range("E6").Value = "'12/10/2022"
$ws2.Range("E6").NumberFormat = "mm-dd-yy"

$ws2.Range("E9").Value = "'12/17/2022"
$ws2.Range("E9").NumberFormat = "mm-dd-yy"

# --- Sheet1: update release dates (column F) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F3").Value = "'10/8/2022"
$ws1.Range("F3").NumberFormat = "mm-dd-yy"

$ws1.Range("F6").Value = "'10/15/2022"
$ws1.Range("F6").NumberFormat = "mm-dd-yy"

$ws1.Range("F8").Value = "'10/29/2022"
$ws1.Range("F8").NumberFormat = "mm-dd-yy"

$ws1.Range("F13").Value = "'11/12/2022"
$ws1.Range("F13").NumberFormat = "mm-dd-yy"

$ws1.Range("F11").Value = "'10/4/2022"
$ws1.Range("F11").NumberFormat = "mm-dd-yy"

$ws1.Range("F15").Value = "'11/26/2022"
$ws1.Range("F15").NumberFormat = "mm-dd-yy"
